$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO"
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M4").Value = 54.26
$wsGrupo.Range("P4").Value = 2.12
$wsGrupo.Range("M10").Value = "1 de 8"
$wsGrupo.Range("P10").Value = "1 de 8"

# Sheet "VENTA MENSUAL"
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F4").Value = 56.38
$wsMensual.Range("F10").Value = 56.38
